# 2/9/25 falta solo postergar semanas y tropas
#
# Mark the F:I "x" columns for rows 47, 49, 51-60 with a yellow highlight
# (reusing the workbook's existing highlighted "x" style already used for
# rows above, e.g. row 43/45/46), and flag rows 59-60 with an upper-case
# "X" instead of lower-case "x". Also move the frozen-pane viewport /
# active selection down to where the user left off (F61).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows that get a lower-case "x"
$xRows = 47, 49, 51, 52, 53, 54, 55, 56, 57, 58

foreach ($r in $xRows) {
    foreach ($col in "F", "G", "H", "I") {
        $cell = $ws.Range("$col$r")
        $cell.Interior.Color = 65535
        $cell.Value = "x"
    }
}

# Rows that get an upper-case "X"
$XRows = 59, 60

foreach ($r in $XRows) {
    foreach ($col in "F", "G", "H", "I") {
        $cell = $ws.Range("$col$r")
        $cell.Interior.Color = 65535
        $cell.Value = "X"
    }
}

# Scroll the frozen pane down and move the active selection to F61
$excel.ActiveWindow.ScrollRow = 47
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F61").Select()
